# Businesspartner, employee bulk addition validations
# - Insert a new "Fax" column (Y) before "Email Address", shifting subsequent
#   columns one to the right (Y..AF -> Z..AG).
# - Populate the new Fax header/value and refresh the sample data row.
# - Scroll the sheet view so column R is the top-left visible column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at Y (25) -- everything from Y onward shifts right.
$ws.Columns.Item(25).EntireColumn.Insert()

# Header row (row 1)
$ws.Cells.Item(1, 25).Value = "Fax"

# Data row (row 2) -- full refresh of the sample record per the new layout.
# (Fax value set first so new shared strings land in the same order as the
# authored workbook: Fax, test fax, master id check, web camera, sales.)
$ws.Cells.Item(2, 25).Value = "test fax"
$ws.Cells.Item(2, 1).Value  = "master id check"
$ws.Cells.Item(2, 2).Value  = "web camera"
$ws.Cells.Item(2, 3).Value  = "sales"
$ws.Cells.Item(2, 4).Value  = "india"
$ws.Cells.Item(2, 5).Value  = "tamilnadu"
$ws.Cells.Item(2, 6).Value  = "madurai"
$ws.Cells.Item(2, 7).Value  = "madurai"
$ws.Cells.Item(2, 8).Value  = "tamilnadu"
$ws.Cells.Item(2, 9).Value  = "india"
$ws.Cells.Item(2, 10).Value = 625144
$ws.Cells.Item(2, 11).Value = "india"
$ws.Cells.Item(2, 12).Value = "tamilnadu"
$ws.Cells.Item(2, 13).Value = "madurai"
$ws.Cells.Item(2, 14).Value = "madurai"
$ws.Cells.Item(2, 15).Value = "tamilnadu"
$ws.Cells.Item(2, 16).Value = "india"
$ws.Cells.Item(2, 17).Value = 625020
$ws.Cells.Item(2, 18).Value = $true
$ws.Cells.Item(2, 19).Value = "Camera"
$ws.Cells.Item(2, 20).Value = "Troy"
$ws.Cells.Item(2, 21).Value = "KKThoppu"
$ws.Cells.Item(2, 22).Value = 2653333
$ws.Cells.Item(2, 23).Value = 2598888
$ws.Cells.Item(2, 24).Value = 9955869555
$ws.Cells.Item(2, 26).Value = "aaa@troy-plus.co.in"
$ws.Cells.Item(2, 27).Value = "www.troy-plus.co.in"
$ws.Cells.Item(2, 28).Value = "hr"
$ws.Cells.Item(2, 29).Value = "test remarks"
$ws.Cells.Item(2, 30).Value = "asd"
$ws.Cells.Item(2, 31).Value = "Sales"
$ws.Cells.Item(2, 32).Value = 250000
$ws.Cells.Item(2, 33).Value = 41907

# Keep the view scrolled like the saved workbook (topLeftCell="R1").
$ws.Application.ActiveWindow.ScrollColumn = 18
